$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 448.1111
$ws.Range("J9").Value = 497.57144
$ws.Range("L9").Value = 497.57144
$ws.Range("N9").Value = -835.5714399999999
$ws.Range("H107").Value = 1161.3889
$ws.Range("I107").Value = 1134.6923
$ws.Range("J107").Value = 1230.8
$ws.Range("K107").Value = 1134.6923
$ws.Range("L107").Value = 1230.8
$ws.Range("M107").Value = 785.3077000000001
$ws.Range("N107").Value = -5070.8
$ws.Range("H132").Value = 13967.391
$ws.Range("I132").Value = 1350.3334
$ws.Range("K132").Value = 4051.0002
$ws.Range("M132").Value = -1521.0002
$ws.Range("H134").Value = 38952.383
$ws.Range("J134").Value = 38952.383
$ws.Range("L134").Value = 38952.383
$ws.Range("N134").Value = -49092.383
$ws.Range("H135").Value = 12620.091
$ws.Range("I135").Value = 3364.2
$ws.Range("J135").Value = 20333.334
$ws.Range("K135").Value = 30277.8
$ws.Range("L135").Value = 183000.006
$ws.Range("M135").Value = -27742.8
$ws.Range("N135").Value = -188070.006
$ws.Range("H136").Value = 40000
$ws.Range("J136").Value = 40000
$ws.Range("L136").Value = 40000
$ws.Range("N136").Value = -50200
$ws.Range("H139").Value = 45682.855
$ws.Range("J139").Value = 45682.855
$ws.Range("L139").Value = 45682.855
$ws.Range("N139").Value = -55962.855
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 14766.5
$ws.Range("J2").Value = 4678.778
$ws.Range("L2").Value = 4678.778
$ws.Range("N2").Value = -4904.778
$ws.Range("H116").Value = 14766.5
$ws.Range("J116").Value = 4678.778
$ws.Range("L116").Value = 4678.778
$ws.Range("N116").Value = -9266.778
$ws.Range("H122").Value = 5799.76
$ws.Range("I122").Value = 5642.4287
$ws.Range("K122").Value = 16927.2861
$ws.Range("M122").Value = -14477.2861

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 14766.5
$ws.Range("J3").Value = 4678.778
$ws.Range("L3").Value = 4678.778
$ws.Range("N3").Value = -4906.778
$ws.Range("H20").Value = 4350.7
$ws.Range("I20").Value = 2501.1428
$ws.Range("J20").Value = 8666.333000000001
$ws.Range("K20").Value = 2501.1428
$ws.Range("L20").Value = 8666.333000000001
$ws.Range("M20").Value = -2254.1428
$ws.Range("N20").Value = -9160.333000000001
$ws.Range("H86").Value = 25058772
$ws.Range("I86").Value = 27840190
$ws.Range("K86").Value = 27840190
$ws.Range("M86").Value = -27839067
$ws.Range("H89").Value = 25058772
$ws.Range("I89").Value = 27840190
$ws.Range("K89").Value = 139200950
$ws.Range("M89").Value = -139195334
$ws.Range("H134").Value = 1236.3889
$ws.Range("I134").Value = 1015.05884
$ws.Range("K134").Value = 3045.17652
$ws.Range("M134").Value = -510.17652

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1366.25
$ws.Range("I16").Value = 1166.1111
$ws.Range("K16").Value = 1166.1111
$ws.Range("M16").Value = -879.1111000000001
$ws.Range("H22").Value = 365.88235
$ws.Range("I22").Value = 426
$ws.Range("K22").Value = 426
$ws.Range("M22").Value = -76
$ws.Range("H86").Value = 26594.4
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 26594.4
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 26594.4
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -28840.4
$ws.Range("H89").Value = 26594.4
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 26594.4
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 132972
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -144204
$ws.Range("H110").Value = 80000
$ws.Range("J110").Value = 80000
$ws.Range("L110").Value = 80000
$ws.Range("N110").Value = -88180
$ws.Range("H113").Value = 1366.25
$ws.Range("I113").Value = 1166.1111
$ws.Range("K113").Value = 1166.1111
$ws.Range("M113").Value = 1003.8889
$ws.Range("H132").Value = 1944.2222
$ws.Range("I132").Value = 2047.7667
$ws.Range("J132").Value = 1426.5
$ws.Range("K132").Value = 6143.300099999999
$ws.Range("L132").Value = 4279.5
$ws.Range("M132").Value = -3613.300099999999
$ws.Range("N132").Value = -9339.5
$ws.Range("H138").Value = 52854.5
$ws.Range("J138").Value = 50000
$ws.Range("L138").Value = 50000
$ws.Range("N138").Value = -60280

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 162990.75
$ws.Range("I70").Value = 1250
$ws.Range("K70").Value = 3750
$ws.Range("M70").Value = -3435
$ws.Range("H73").Value = 162990.75
$ws.Range("I73").Value = 1250
$ws.Range("K73").Value = 3750
$ws.Range("M73").Value = -2658
$ws.Range("H131").Value = 5930
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 5930
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 17790
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -27870

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 58562.81
$ws.Range("I70").Value = 79471.47
$ws.Range("J70").Value = 6291.1665
$ws.Range("K70").Value = 79471.47
$ws.Range("L70").Value = 6291.1665
$ws.Range("M70").Value = -79201.47
$ws.Range("N70").Value = -6831.1665
$ws.Range("H73").Value = 58562.81
$ws.Range("I73").Value = 79471.47
$ws.Range("J73").Value = 6291.1665
$ws.Range("K73").Value = 79471.47
$ws.Range("L73").Value = 6291.1665
$ws.Range("M73").Value = -78535.47
$ws.Range("N73").Value = -8163.1665
$ws.Range("H94").Value = 28263.273
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 28263.273
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 28263.273
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -29615.273
$ws.Range("H132").Value = 6975.213
$ws.Range("I132").Value = 6282.5264
$ws.Range("K132").Value = 18847.5792
$ws.Range("M132").Value = -16317.5792

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3418.4
$ws.Range("I46").Value = 1500
$ws.Range("K46").Value = 1500
$ws.Range("M46").Value = -1312
$ws.Range("H55").Value = 402.44446
$ws.Range("I55").Value = 490.16666
$ws.Range("K55").Value = 490.16666
$ws.Range("M55").Value = -317.16666
$ws.Range("H61").Value = 1327.3
$ws.Range("I61").Value = 1327.3
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1327.3
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1125.3
$ws.Range("N61").ClearContents()
$ws.Range("H113").Value = 1327.3
$ws.Range("I113").Value = 1327.3
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1327.3
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 842.7
$ws.Range("N113").ClearContents()
$ws.Range("H136").Value = 3425.375
$ws.Range("I136").Value = 4333.25
$ws.Range("J136").Value = 2517.5
$ws.Range("K136").Value = 12999.75
$ws.Range("L136").Value = 7552.5
$ws.Range("M136").Value = -10449.75
$ws.Range("N136").Value = -12652.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 5842.3335
$ws.Range("I26").Value = 7097
$ws.Range("K26").Value = 7097
$ws.Range("M26").Value = -6804
$ws.Range("H98").Value = 49795
$ws.Range("J98").Value = 49795
$ws.Range("L98").Value = 49795
$ws.Range("N98").Value = -55785
